$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#convert")
$ws.Activate()

# The "#unique=false" tag was renamed to "#match=all" throughout the
# testing files. Update every cell in column D that still holds the old
# value.
$rows = @(1, 4, 8, 11, 15, 18)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq "#unique=false") {
        $cell.Value2 = "#match=all"
    }
}

$ws.Range("F12").Select()
